# Update the sensor X/Y coordinates (columns E "x_m" and F "y_m", rows
# 2-49) on the active worksheet to reflect the corrected projection used
# after flipping the sensor-plan image. Column G ("height_m") is left
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value2 = 23.69979421125619
$ws.Cells.Item(2, 6).Value2 = 1.767946064993388
$ws.Cells.Item(3, 5).Value2 = 0.2821262680420327
$ws.Cells.Item(3, 6).Value2 = 1.837069358881833
$ws.Cells.Item(4, 5).Value2 = 0.2821262680420327
$ws.Cells.Item(4, 6).Value2 = 1.837069358881833
$ws.Cells.Item(5, 5).Value2 = 0.2821262680420327
$ws.Cells.Item(5, 6).Value2 = 1.837069358881833
$ws.Cells.Item(6, 5).Value2 = 0.2821262680420327
$ws.Cells.Item(6, 6).Value2 = 1.837069358881833
$ws.Cells.Item(7, 5).Value2 = 22.7329790285772
$ws.Cells.Item(7, 6).Value2 = 1.74260225100966
$ws.Cells.Item(8, 5).Value2 = 5.74268778331578
$ws.Cells.Item(8, 6).Value2 = 2.801719087524943
$ws.Cells.Item(9, 5).Value2 = 5.74268778331578
$ws.Cells.Item(9, 6).Value2 = 2.801719087524943
$ws.Cells.Item(10, 5).Value2 = 5.74268778331578
$ws.Cells.Item(10, 6).Value2 = 2.801719087524943
$ws.Cells.Item(11, 5).Value2 = 5.74268778331578
$ws.Cells.Item(11, 6).Value2 = 2.801719087524943
$ws.Cells.Item(12, 5).Value2 = 8.076524319340663
$ws.Cells.Item(12, 6).Value2 = 1.658488031914387
$ws.Cells.Item(13, 5).Value2 = 8.076524319340663
$ws.Cells.Item(13, 6).Value2 = 1.658488031914387
$ws.Cells.Item(14, 5).Value2 = 8.076524319340663
$ws.Cells.Item(14, 6).Value2 = 1.658488031914387
$ws.Cells.Item(15, 5).Value2 = 8.076524319340663
$ws.Cells.Item(15, 6).Value2 = 1.658488031914387
$ws.Cells.Item(16, 5).Value2 = 6.909606051328223
$ws.Cells.Item(16, 6).Value2 = 2.230103559719665
$ws.Cells.Item(17, 5).Value2 = 6.909606051328223
$ws.Cells.Item(17, 6).Value2 = 2.230103559719665
$ws.Cells.Item(18, 5).Value2 = 6.909606051328223
$ws.Cells.Item(18, 6).Value2 = 2.230103559719665
$ws.Cells.Item(19, 5).Value2 = 6.909606051328223
$ws.Cells.Item(19, 6).Value2 = 2.230103559719665
$ws.Cells.Item(20, 5).Value2 = 6.909606051328223
$ws.Cells.Item(20, 6).Value2 = 2.230103559719665
$ws.Cells.Item(21, 5).Value2 = 6.909606051328223
$ws.Cells.Item(21, 6).Value2 = 2.230103559719665
$ws.Cells.Item(22, 5).Value2 = 10.41036085536555
$ws.Cells.Item(22, 6).Value2 = 2.801719087524943
$ws.Cells.Item(23, 5).Value2 = 10.41036085536555
$ws.Cells.Item(23, 6).Value2 = 2.801719087524943
$ws.Cells.Item(24, 5).Value2 = 10.41036085536555
$ws.Cells.Item(24, 6).Value2 = 2.801719087524943
$ws.Cells.Item(25, 5).Value2 = 10.41036085536555
$ws.Cells.Item(25, 6).Value2 = 2.801719087524943
$ws.Cells.Item(26, 5).Value2 = 12.74419739139043
$ws.Cells.Item(26, 6).Value2 = 1.658488031914387
$ws.Cells.Item(27, 5).Value2 = 12.74419739139043
$ws.Cells.Item(27, 6).Value2 = 1.658488031914387
$ws.Cells.Item(28, 5).Value2 = 12.74419739139043
$ws.Cells.Item(28, 6).Value2 = 1.658488031914387
$ws.Cells.Item(29, 5).Value2 = 12.74419739139043
$ws.Cells.Item(29, 6).Value2 = 1.658488031914387
$ws.Cells.Item(30, 5).Value2 = 11.57727912337799
$ws.Cells.Item(30, 6).Value2 = 2.230103559719665
$ws.Cells.Item(31, 5).Value2 = 11.57727912337799
$ws.Cells.Item(31, 6).Value2 = 2.230103559719665
$ws.Cells.Item(32, 5).Value2 = 11.57727912337799
$ws.Cells.Item(32, 6).Value2 = 2.230103559719665
$ws.Cells.Item(33, 5).Value2 = 11.57727912337799
$ws.Cells.Item(33, 6).Value2 = 2.230103559719665
$ws.Cells.Item(34, 5).Value2 = 11.57727912337799
$ws.Cells.Item(34, 6).Value2 = 2.230103559719665
$ws.Cells.Item(35, 5).Value2 = 11.57727912337799
$ws.Cells.Item(35, 6).Value2 = 2.230103559719665
$ws.Cells.Item(36, 5).Value2 = 1.563964053029377
$ws.Cells.Item(36, 6).Value2 = 2.692495134041453
$ws.Cells.Item(37, 5).Value2 = 1.563964053029377
$ws.Cells.Item(37, 6).Value2 = 2.692495134041453
$ws.Cells.Item(38, 5).Value2 = 1.563964053029377
$ws.Cells.Item(38, 6).Value2 = 2.692495134041453
$ws.Cells.Item(39, 5).Value2 = 1.563964053029377
$ws.Cells.Item(39, 6).Value2 = 2.692495134041453
$ws.Cells.Item(40, 5).Value2 = 3.562529915067857
$ws.Cells.Item(40, 6).Value2 = 0.9051969958217355
$ws.Cells.Item(41, 5).Value2 = 3.562529915067857
$ws.Cells.Item(41, 6).Value2 = 0.9051969958217355
$ws.Cells.Item(42, 5).Value2 = 3.562529915067857
$ws.Cells.Item(42, 6).Value2 = 0.9051969958217355
$ws.Cells.Item(43, 5).Value2 = 3.562529915067857
$ws.Cells.Item(43, 6).Value2 = 0.9051969958217355
$ws.Cells.Item(44, 5).Value2 = 2.563246984048617
$ws.Cells.Item(44, 6).Value2 = 1.798846064931594
$ws.Cells.Item(45, 5).Value2 = 2.563246984048617
$ws.Cells.Item(45, 6).Value2 = 1.798846064931594
$ws.Cells.Item(46, 5).Value2 = 2.563246984048617
$ws.Cells.Item(46, 6).Value2 = 1.798846064931594
$ws.Cells.Item(47, 5).Value2 = 2.563246984048617
$ws.Cells.Item(47, 6).Value2 = 1.798846064931594
$ws.Cells.Item(48, 5).Value2 = 2.563246984048617
$ws.Cells.Item(48, 6).Value2 = 1.798846064931594
$ws.Cells.Item(49, 5).Value2 = 2.563246984048617
$ws.Cells.Item(49, 6).Value2 = 1.798846064931594
